$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump version & date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Elements sheet: add ele-1/ext-1 invariants text to the root Extension row (row 1, column AJ) ---
$elems = $wb.Worksheets.Item("Elements")
$elems.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
